# Natmi following Dr Hou advice
# Expands the LR-pair table from a single sending/target cluster pair (FAPs -> FAPs, sCs)
# to the full combination of sending clusters {ECs, FAPs} x target clusters {ECs, FAPs, sCs},
# with recomputed specificity / weight statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Wnt1", "Ror2", "ECs",  1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 1, 0.3333333333333333, 0.1350926666666667, 0.405278,    0.03398966953169371, 0.03398966953169372, 0.003664118398,     0.032977065582,     0.002423954499232144, 0.002423954499232145),
    @("ECs",  "Wnt1", "Ror2", "FAPs", 1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1,                  2.244570333333333, 6.733711,    0.5647397875333248,  0.5647397875333248,  0.060879481151,     0.547915330359,     0.04027410586061662,  0.04027410586061662),
    @("ECs",  "Wnt1", "Ror2", "sCs",  1, 0.3333333333333333, 0.027123, 0.081369, 0.07131444737854614, 0.07131444737854616, 3, 1,                  1.594858333333333, 4.784575,    0.4012705429349815,  0.4012705429349816,  0.043257342575,     0.389316083175,     0.02861638701869738,  0.02861638701869739),
    @("FAPs", "Wnt1", "Ror2", "ECs",  3, 1,                  0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 1, 0.3333333333333333, 0.1350926666666667, 0.405278,    0.03398966953169371, 0.03398966953169372, 0.04771563048444444, 0.4294406743600001, 0.03156571503246156,  0.03156571503246157),
    @("FAPs", "Wnt1", "Ror2", "FAPs", 3, 1,                  0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1,                  2.244570333333333, 6.733711,    0.5647397875333248,  0.5647397875333248,  0.7927972055355555, 7.135174849819999,  0.5244656816727081,   0.5244656816727081),
    @("FAPs", "Wnt1", "Ror2", "sCs",  3, 1,                  0.3532066666666667, 1.05962, 0.9286855526214538, 0.9286855526214538, 3, 1,                  1.594858333333333, 4.784575,    0.4012705429349815,  0.4012705429349816,  0.5633145957222223, 5.0698313615,       0.3726541559162841,   0.3726541559162842)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($c = 0; $c -lt $values.Count; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $values[$c]
    }
}
